# feat: add 2022-Q4 data
#
# 1) Insert a brand-new "2022-Q4" sheet (fund-holdings detail) right after
#    "总计" and before "2022-Q3" - built by duplicating the "2022-Q3" sheet
#    (so formatting / text-vs-number cell typing matches the existing
#    quarter sheets) and then overwriting the handful of cells whose values
#    actually differ, trimming the two rows that don't apply to this quarter.
# 2) Update the "总计" (summary) sheet so the new quarter becomes the first
#    data row and every later row's date/count/value shifts down by one,
#    with a brand new trailing row for what used to be the oldest (now
#    second-oldest) quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the "2022-Q4" worksheet
# ---------------------------------------------------------------------
$srcQ3 = $wb.Worksheets.Item("2022-Q3")
$srcQ3.Copy($srcQ3, $null)          # places the copy immediately before 2022-Q3
$newQ4 = $wb.Worksheets.Item("2022-Q3 (2)")
$newQ4.Name = "2022-Q4"

# 2022-Q3 had 4 funds, 2022-Q4 only has 2 - drop the extra two rows.
$newQ4.Rows("4:5").Delete()

# Row 2 - 517160 / 南方中证长江保护主题ETF (code/name/scale stay "517160" text,
# only scale/position/weight/value/rank differ from the Q3 sheet we copied).
$newQ4.Range("D2").NumberFormat = "@"
$newQ4.Range("D2").Value = "16.77"
$newQ4.Range("E2").NumberFormat = "@"
$newQ4.Range("E2").Value = "99.26"
$newQ4.Range("F2").NumberFormat = "@"
$newQ4.Range("F2").Value = "2.41"
$newQ4.Range("G2").NumberFormat = "@"
$newQ4.Range("G2").Value = "0.4042"
$newQ4.Range("H2").Value = 7

# Row 3 - 517330 / 易方达中证长江保护主题ETF
$newQ4.Range("D3").NumberFormat = "@"
$newQ4.Range("D3").Value = "16.47"
$newQ4.Range("E3").NumberFormat = "@"
$newQ4.Range("E3").Value = "99.39"
$newQ4.Range("F3").NumberFormat = "@"
$newQ4.Range("F3").Value = "2.41"
$newQ4.Range("G3").NumberFormat = "@"
$newQ4.Range("G3").Value = "0.3969"
$newQ4.Range("H3").Value = 7

# ---------------------------------------------------------------------
# Step 2: refresh the "总计" summary table
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")

$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = 0.8

$ws1.Range("B3").Value = "2022-Q3"
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 0.83

$ws1.Range("B4").Value = "2022-Q1"
$ws1.Range("C4").Value = 3
$ws1.Range("D4").Value = 4.81

$ws1.Range("B5").Value = "2021-Q4"
$ws1.Range("C5").Value = 5
$ws1.Range("D5").Value = 8.15

$ws1.Range("B6").Value = "2021-Q3"
$ws1.Range("C6").Value = 11
$ws1.Range("D6").Value = 10.55

$ws1.Range("B7").Value = "2021-Q2"
$ws1.Range("C7").Value = 5
$ws1.Range("D7").Value = 3.61

$ws1.Range("B8").Value = "2021-Q1"
$ws1.Range("C8").Value = 2
$ws1.Range("D8").Value = 0.51

# New trailing row for 2020-Q4 (row 9). Copy A8's number formatting/style
# onto A9 before writing its value, so it matches the other index cells.
$ws1.Range("A8").Copy()
$ws1.Range("A9").PasteSpecial(-4122)
$ws1.Range("A9").Value = 7
$ws1.Range("B9").Value = "2020-Q4"
$ws1.Range("C9").Value = 3
$ws1.Range("D9").Value = 0.54
